# Add two new daily rows (2025-11-24 and 2025-11-25, serials 45985/45986)
# to the bottom of each of the 8 stock sheets, extending the used range
# from A1:B115 to A1:B117.

$wb = $excel.ActiveWorkbook

# New data per sheet (by tab order / index), matching column B (remn_amt)
# values added in the upstream diff.
$newValues = @{
    1 = @(10150878, 10302767)  # 삼성전자
    2 = @(11050879, 10649732)  # SK하이닉스
    3 = @(2897779,  2883157)   # LG에너지솔루션
    4 = @(834399,   844863)    # LG화학
    5 = @(1401580,  1388186)   # 에코프로비엠
    6 = @(1623538,  1617790)   # 에코프로
    7 = @(248493,   251432)    # LG전자
    8 = @(287046,   286145)    # LG디스플레이
}

$dateSerials = @(45985, 45986)

for ($i = 1; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $values = $newValues[$i]

    for ($j = 0; $j -lt 2; $j++) {
        $row = 116 + $j
        $ws.Cells.Item($row, 1).Value = $dateSerials[$j]
        $ws.Cells.Item($row, 2).Value = $values[$j]
    }

    # Match the existing date-column formatting used throughout column A.
    $ws.Range("A116:A117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
